$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "Save" column in H1, using the same style as the
# other header cells (style index 1 -> bold, centered, bordered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the Save values for rows 2-9
$saveValues = @(0, 0, 1, 0, 0, 1, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
